# Auto-generated edit script: updates market-price derived columns (H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# refreshed Universalis market data, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K18").Value = 2433.3333
$ws.Range("H18").Value = 4325
$ws.Range("I18").Value = 2433.3333
$ws.Range("M18").Value = -2149.3333
$ws.Range("I31").Value = 2249.75
$ws.Range("H31").Value = 2249.75
$ws.Range("K31").Value = 6749.25
$ws.Range("M31").Value = -6519.25
$ws.Range("H116").Value = 4951.25
$ws.Range("I116").Value = 4902.5
$ws.Range("K116").Value = 4902.5
$ws.Range("J116").Value = 5000
$ws.Range("M116").Value = -1460.5
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H132").Value = 1903.3167
$ws.Range("M132").Value = -3291.6314
$ws.Range("I132").Value = 1940.5438
$ws.Range("K132").Value = 5821.6314
$ws.Range("J133").Value = 94999.8
$ws.Range("N133").Value = -105119.8
$ws.Range("H133").Value = 94999.8
$ws.Range("L133").Value = 94999.8
$ws.Range("K141").Value = 2881.7586
$ws.Range("M141").Value = 2298.2414
$ws.Range("H141").Value = 1153.4517
$ws.Range("I141").Value = 960.5862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I6").Value = 5000
$ws.Range("H6").Value = 5000
$ws.Range("M6").Value = -4827
$ws.Range("K6").Value = 5000
$ws.Range("M46").ClearContents()
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 12766.5
$ws.Range("N46").Value = -13404.5
$ws.Range("H46").Value = 12766.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 12766.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M22").Value = -530.25
$ws.Range("H22").Value = 357791.16
$ws.Range("I22").Value = 703.25
$ws.Range("K22").Value = 703.25
$ws.Range("I86").Value = 6751.3335
$ws.Range("M86").Value = -5628.3335
$ws.Range("K86").Value = 6751.3335
$ws.Range("H86").Value = 5736.684
$ws.Range("I89").Value = 6751.3335
$ws.Range("M89").Value = -28140.6675
$ws.Range("K89").Value = 33756.6675
$ws.Range("H89").Value = 5736.684
$ws.Range("H95").Value = 12905.75
$ws.Range("N95").Value = -18397.75
$ws.Range("L95").Value = 12905.75
$ws.Range("J95").Value = 12905.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M22").Value = -18.19999999999999
$ws.Range("H22").Value = 640
$ws.Range("I22").Value = 368.2
$ws.Range("K22").Value = 368.2
$ws.Range("I31").Value = 1564.4445
$ws.Range("H31").Value = 1733.6333
$ws.Range("K31").Value = 1564.4445
$ws.Range("M31").Value = -1269.4445
$ws.Range("K34").Value = 1564.4445
$ws.Range("H34").Value = 1733.6333
$ws.Range("I34").Value = 1564.4445
$ws.Range("M34").Value = -1362.4445
$ws.Range("M134").Value = -38986.5
$ws.Range("I134").Value = 13840.5
$ws.Range("H134").Value = 17015.912
$ws.Range("K134").Value = 41521.5
$ws.Range("N141").Value = -229189.8
$ws.Range("L141").Value = 218829.8
$ws.Range("H141").Value = 218829.8
$ws.Range("J141").Value = 218829.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 3607078.5
$ws.Range("N4").Value = -7093147.399999999
$ws.Range("K4").Value = 10821235.5
$ws.Range("J4").Value = 2364307.8
$ws.Range("L4").Value = 7092923.399999999
$ws.Range("M4").Value = -10821123.5
$ws.Range("H4").Value = 3353921.5
$ws.Range("I112").Value = 979.5
$ws.Range("J112").Value = 9999
$ws.Range("K112").Value = 2938.5
$ws.Range("L112").Value = 29997
$ws.Range("M112").Value = -1830.5
$ws.Range("H112").Value = 6391.2
$ws.Range("N112").Value = -32213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20926
$ws.Range("H38").Value = 20000
$ws.Range("J38").Value = 20000
$ws.Range("L39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("J39").Value = 0
$ws.Range("K97").Value = 2780.625
$ws.Range("M97").Value = -2284.625
$ws.Range("N97").Value = -2356
$ws.Range("L97").Value = 1364
$ws.Range("J97").Value = 1364
$ws.Range("I97").Value = 2780.625
$ws.Range("H97").Value = 2173.5
$ws.Range("J139").Value = 65000
$ws.Range("H139").Value = 73432
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280
$ws.Range("J140").Value = 113744.75
$ws.Range("H140").Value = 113744.75
$ws.Range("N140").Value = -124104.75
$ws.Range("L140").Value = 113744.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 5950.2
$ws.Range("I7").Value = 5950.2
$ws.Range("H7").Value = 6718.625
$ws.Range("M7").Value = -5838.2
$ws.Range("H20").Value = 12560.2
$ws.Range("J20").Value = 10933.667
$ws.Range("N20").Value = -11385.667
$ws.Range("L20").Value = 10933.667
$ws.Range("M122").Value = -18548.9995
$ws.Range("K122").Value = 20998.9995
$ws.Range("H122").Value = 7414.7
$ws.Range("I122").Value = 6999.6665
$ws.Range("I126").Value = 5950.2
$ws.Range("K126").Value = 17850.6
$ws.Range("M126").Value = -15380.6
$ws.Range("H126").Value = 6718.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J32").Value = 24495
$ws.Range("N32").Value = -25129
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 24495
$ws.Range("M32").ClearContents()
$ws.Range("H32").Value = 24495
$ws.Range("K34").Value = 505000
$ws.Range("H34").Value = 505000
$ws.Range("I34").Value = 505000
$ws.Range("M34").Value = -504797
$ws.Range("N81").Value = -12121.6474
$ws.Range("I81").Value = 9398.714
$ws.Range("J81").Value = 4999.8237
$ws.Range("L81").Value = 9999.6474
$ws.Range("K81").Value = 18797.428
$ws.Range("H81").Value = 6986.4194
$ws.Range("M81").Value = -17736.428
$ws.Range("L84").Value = 49998.237
$ws.Range("K84").Value = 93987.14
$ws.Range("M84").Value = -88683.14
$ws.Range("H84").Value = 6986.4194
$ws.Range("J84").Value = 4999.8237
$ws.Range("I84").Value = 9398.714
$ws.Range("N84").Value = -60606.237
$ws.Range("L94").Value = 28083.334
$ws.Range("H94").Value = 27642.857
$ws.Range("N94").Value = -29885.334
$ws.Range("J94").Value = 28083.334
$ws.Range("I139").Value = 5000
$ws.Range("J139").Value = 63637.273
$ws.Range("H139").Value = 58750.832
$ws.Range("M139").Value = 140
$ws.Range("L139").Value = 63637.273
$ws.Range("N139").Value = -73917.273
$ws.Range("K139").Value = 5000
$ws.Range("N141").Value = -131608.75
$ws.Range("K141").Value = 74000
$ws.Range("L141").Value = 121248.75
$ws.Range("M141").Value = -68820
$ws.Range("H141").Value = 105499.164
$ws.Range("I141").Value = 74000
$ws.Range("J141").Value = 121248.75
